$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-12 correspond to N = 1..11 (row - 1)
# Column C (nama): "Tanah PB N" -> "Tanah Lapangan N"
# Column F (letak_tanah): "Jalan Prabubima N" -> "Jalan Gerilya N"
# Column J (penggunaan): "Kantor PB N" -> "Tempat Olahraga N"
# Column L (keterangan): "Milik Pak Aribawa" -> "Milik Pemda BMS"
for ($row = 2; $row -le 12; $row++) {
    $n = $row - 1
    $ws.Cells.Item($row, 3).Value = "Tanah Lapangan $n"
    $ws.Cells.Item($row, 6).Value = "Jalan Gerilya $n"
    $ws.Cells.Item($row, 10).Value = "Tempat Olahraga $n"
    $ws.Cells.Item($row, 12).Value = "Milik Pemda BMS"
}

# Column width adjustments (nama / penggunaan columns widened after content change)
$ws.Columns.Item(3).ColumnWidth = 25
$ws.Columns.Item(10).ColumnWidth = 16.83

# Update the active selection to reflect where the user ended up
[void]$ws.Range("F15").Select()

Write-Host "Done"
